$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.378.35"
$ws.Range("E2").Value = "  -2.21%  "
$ws.Range("D3").Value = "1.709.15"
$ws.Range("E3").Value = "  -1.82%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("E5").Value = "  -1.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5316"
$ws.Range("E6").Value = "  -1.95%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.005"
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  -3.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06587"
$ws.Range("E9").Value = "  -1.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.81"
$ws.Range("E10").Value = "  -3.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07634"
$ws.Range("E11").Value = "  -1.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.572"
$ws.Range("E12").Value = "  -2.59%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.719.57"
$ws.Range("E13").Value = "  -1.60%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "1.945.97"
$ws.Range("E14").Value = "  -1.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5722"
$ws.Range("E15").Value = "  -3.93%  "
$ws.Range("D16").Value = "0.0₅8163"
$ws.Range("E16").Value = "  -2.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.87"
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").Value = "27.356.29"
$ws.Range("E18").Value = "  -2.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.34"
$ws.Range("E19").Value = "  -3.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.005"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.669"
$ws.Range("E21").Value = "  -3.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.43"
$ws.Range("E22").Value = "  -4.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.971"
$ws.Range("E23").Value = "  -4.07%  "
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.771"
$ws.Range("E25").Value = "  +6.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.81"
$ws.Range("E26").Value = "  -2.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1216"
$ws.Range("E27").Value = "  -2.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.272"
$ws.Range("E28").Value = "  -2.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.29"
$ws.Range("E29").Value = "  -5.20%  "
$ws.Range("E30").Value = "  -4.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.292"
$ws.Range("E31").Value = "  -1.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.505"
$ws.Range("E32").Value = "  -5.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.429"
$ws.Range("E33").Value = "  -2.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.645"
$ws.Range("E34").Value = "  -1.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.877"
$ws.Range("E35").Value = "  +0.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9495"
$ws.Range("E36").Value = "  -3.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.414"
$ws.Range("E37").Value = "  -1.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5861"
$ws.Range("E38").Value = "  -1.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01628"
$ws.Range("E39").Value = "  -2.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.865"
$ws.Range("E40").Value = "  -1.07%  "
$ws.Range("D41").Value = "1.046.27"
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.005"
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8441"
$ws.Range("E43").Value = "  -0.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.93"
$ws.Range("E44").Value = "  -1.30%  "
$ws.Range("D45").Value = "1.852.36"
$ws.Range("E45").Value = "  -1.81%  "
$ws.Range("D46").Value = "0.0₈118"
$ws.Range("E46").Value = "  +0.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "58.00"
$ws.Range("E47").Value = "  -3.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4509"
$ws.Range("E48").Value = "  +1.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.004"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.080"
$ws.Range("E50").Value = "  -2.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05245"
$ws.Range("E51").Value = "  -1.44%  "
